$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (ALC)
$ws.Range("H33").Value = 332.70834
$ws.Range("I33").Value = 353.77274
$ws.Range("K33").Value = 353.77274
$ws.Range("M33").Value = -124.77274

# Row 49 (ALC)
$ws.Range("H49").Value = 487.1111
$ws.Range("I49").Value = 487.1111
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 1461.3333
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1325.3333
$ws.Range("N49").ClearContents()

# Row 98 (ALC)
$ws.Range("H98").Value = 1728.4
$ws.Range("I98").Value = 1823
$ws.Range("J98").Value = 1350
$ws.Range("K98").Value = 1823
$ws.Range("L98").Value = 1350
$ws.Range("M98").Value = -325
$ws.Range("N98").Value = -4346

# Row 107 (ALC)
$ws.Range("H107").Value = 6116.8
$ws.Range("I107").Value = 6116.8
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 6116.8
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -4196.8
$ws.Range("N107").ClearContents()

# Row 122 (ALC)
$ws.Range("H122").Value = 1728.4
$ws.Range("I122").Value = 1823
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 5469
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -3019
$ws.Range("N122").Value = -8950

# Row 138 (ALC)
$ws.Range("H138").Value = 2217.451
$ws.Range("I138").Value = 1082.9131
$ws.Range("J138").Value = 3149.3928
$ws.Range("K138").Value = 3248.7393
$ws.Range("L138").Value = 9448.178400000001
$ws.Range("M138").Value = 1891.2607
$ws.Range("N138").Value = -19728.1784

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (BSM)
$ws.Range("H86").Value = 27780508
$ws.Range("I86").Value = 2717
$ws.Range("K86").Value = 2717
$ws.Range("M86").Value = -1594

# Row 89 (BSM)
$ws.Range("H89").Value = 27780508
$ws.Range("I89").Value = 2717
$ws.Range("K89").Value = 13585
$ws.Range("M89").Value = -7969

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 5955509.5
$ws.Range("I31").Value = 3191.5095
$ws.Range("K31").Value = 3191.5095
$ws.Range("M31").Value = -2896.5095

# Row 34 (CRP)
$ws.Range("H34").Value = 5955509.5
$ws.Range("I34").Value = 3191.5095
$ws.Range("K34").Value = 3191.5095
$ws.Range("M34").Value = -2989.5095

# Row 58 (CRP)
$ws.Range("H58").Value = 3087.3572
$ws.Range("I58").Value = 902.875
$ws.Range("J58").Value = 6000
$ws.Range("K58").Value = 902.875
$ws.Range("L58").Value = 6000
$ws.Range("M58").Value = -699.875
$ws.Range("N58").Value = -6406

# Row 59 (CRP)
$ws.Range("H59").Value = 22877
$ws.Range("J59").Value = 22877
$ws.Range("L59").Value = 22877
$ws.Range("N59").Value = -25167

# Row 105 (CRP)
$ws.Range("H105").Value = 1842.25
$ws.Range("I105").Value = 1955.4286
$ws.Range("J105").Value = 1050
$ws.Range("K105").Value = 1955.4286
$ws.Range("L105").Value = 1050
$ws.Range("M105").Value = -208.4286
$ws.Range("N105").Value = -4544

# Row 107 (CRP)
$ws.Range("H107").Value = 944.86365
$ws.Range("I107").Value = 675.05884
$ws.Range("J107").Value = 1862.2
$ws.Range("K107").Value = 675.05884
$ws.Range("L107").Value = 1862.2
$ws.Range("M107").Value = 1244.94116
$ws.Range("N107").Value = -5702.2

# Row 122 (CRP)
$ws.Range("H122").Value = 1964.4615
$ws.Range("I122").Value = 1967.091
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 5901.272999999999
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3451.272999999999
$ws.Range("N122").Value = -10750

# Row 136 (CRP)
$ws.Range("H136").Value = 3087.3572
$ws.Range("I136").Value = 902.875
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 2708.625
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -158.625
$ws.Range("N136").Value = -23100

# Row 140 (CRP)
$ws.Range("H140").Value = 41375.8
$ws.Range("J140").Value = 41375.8
$ws.Range("L140").Value = 41375.8
$ws.Range("N140").Value = -51735.8

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (CUL)
$ws.Range("H23").Value = 216.875
$ws.Range("J23").Value = 235.66667
$ws.Range("L23").Value = 707.00001
$ws.Range("N23").Value = -1177.00001

# Row 39 (CUL)
$ws.Range("H39").Value = 510.45456
$ws.Range("J39").Value = 510.45456
$ws.Range("L39").Value = 1531.36368
$ws.Range("N39").Value = -2119.36368

# Row 62 (CUL)
$ws.Range("H62").Value = 7937.067
$ws.Range("I62").Value = 7250
$ws.Range("J62").Value = 8186.909
$ws.Range("K62").Value = 21750
$ws.Range("L62").Value = 24560.727
$ws.Range("M62").Value = -21064
$ws.Range("N62").Value = -25932.727

# Row 65 (CUL)
$ws.Range("H65").Value = 7937.067
$ws.Range("I65").Value = 7250
$ws.Range("J65").Value = 8186.909
$ws.Range("K65").Value = 65250
$ws.Range("L65").Value = 73682.181
$ws.Range("M65").Value = -61818
$ws.Range("N65").Value = -80546.181

# Row 87 (CUL)
$ws.Range("H87").Value = 12643.75
$ws.Range("I87").Value = 5200
$ws.Range("J87").Value = 15125
$ws.Range("K87").Value = 15600
$ws.Range("L87").Value = 45375
$ws.Range("M87").Value = -14352
$ws.Range("N87").Value = -47871

# Row 90 (CUL)
$ws.Range("H90").Value = 12643.75
$ws.Range("I90").Value = 5200
$ws.Range("J90").Value = 15125
$ws.Range("K90").Value = 46800
$ws.Range("L90").Value = 136125
$ws.Range("M90").Value = -40560
$ws.Range("N90").Value = -148605

# Row 94 (CUL)
$ws.Range("H94").Value = 4555.5557
$ws.Range("J94").Value = 4500
$ws.Range("L94").Value = 13500
$ws.Range("N94").Value = -14852

# Row 118 (CUL)
$ws.Range("H118").Value = 1103.6
$ws.Range("J118").Value = 1146.421
$ws.Range("L118").Value = 3439.263
$ws.Range("N118").Value = -5925.263

# Row 131 (CUL)
$ws.Range("H131").Value = 814.49
$ws.Range("I131").Value = 356.66666
$ws.Range("J131").Value = 859.7692
$ws.Range("K131").Value = 1069.99998
$ws.Range("L131").Value = 2579.3076
$ws.Range("M131").Value = 3970.00002
$ws.Range("N131").Value = -12659.3076

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (GSM)
$ws.Range("H2").Value = 123.333336
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# Row 102 (GSM)
$ws.Range("H102").Value = 4624.1763
$ws.Range("I102").Value = 4986.5
$ws.Range("J102").Value = 2933.3333
$ws.Range("K102").Value = 4986.5
$ws.Range("L102").Value = 2933.3333
$ws.Range("M102").Value = -3364.5
$ws.Range("N102").Value = -6177.3333

# Row 122 (GSM)
$ws.Range("H122").Value = 8336745
$ws.Range("I122").Value = 16668893
$ws.Range("J122").Value = 4597
$ws.Range("K122").Value = 50006679
$ws.Range("L122").Value = 13791
$ws.Range("M122").Value = -50004229
$ws.Range("N122").Value = -18691

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 1722.75
$ws.Range("I22").Value = 950
$ws.Range("J22").Value = 1877.3
$ws.Range("K22").Value = 950
$ws.Range("L22").Value = 1877.3
$ws.Range("M22").Value = -655
$ws.Range("N22").Value = -2467.3

# Row 27 (LTW)
$ws.Range("H27").Value = 1722.75
$ws.Range("I27").Value = 950
$ws.Range("J27").Value = 1877.3
$ws.Range("K27").Value = 950
$ws.Range("L27").Value = 1877.3
$ws.Range("M27").Value = -843
$ws.Range("N27").Value = -2091.3

# Row 40 (LTW)
$ws.Range("H40").Value = 6813
$ws.Range("J40").Value = 4167.3335
$ws.Range("L40").Value = 4167.3335
$ws.Range("N40").Value = -4439.3335

# Row 61 (LTW)
$ws.Range("H61").Value = 1829.3529
$ws.Range("I61").Value = 1678.2222
$ws.Range("J61").Value = 1999.375
$ws.Range("K61").Value = 1678.2222
$ws.Range("L61").Value = 1999.375
$ws.Range("M61").Value = -1476.2222
$ws.Range("N61").Value = -2403.375

# Row 113 (LTW)
$ws.Range("H113").Value = 1829.3529
$ws.Range("I113").Value = 1678.2222
$ws.Range("J113").Value = 1999.375
$ws.Range("K113").Value = 1678.2222
$ws.Range("L113").Value = 1999.375
$ws.Range("M113").Value = 491.7778000000001
$ws.Range("N113").Value = -6339.375

$ws = $wb.Worksheets.Item("WVR")
# Row 122 (WVR)
$ws.Range("H122").Value = 2870.3
$ws.Range("J122").Value = 1047.5
$ws.Range("L122").Value = 3142.5
$ws.Range("N122").Value = -8042.5

# Row 126 (WVR)
$ws.Range("H126").Value = 2701.5789
$ws.Range("I126").Value = 1816.4286
$ws.Range("J126").Value = 5180
$ws.Range("K126").Value = 5449.2858
$ws.Range("L126").Value = 15540
$ws.Range("M126").Value = -2979.2858
$ws.Range("N126").Value = -20480

# Row 136 (WVR)
$ws.Range("H136").Value = 1292.6316
$ws.Range("I136").Value = 1292.6316
$ws.Range("K136").Value = 3877.8948
$ws.Range("M136").Value = -1327.8948
